# did some leetcode to practice for oa
# - append two new solved problems to the "Arrays & Hashing" tracker sheet
# - add a new "Others" tracker sheet (with the standard header row) at the end

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Append rows 13 & 14 to the "Arrays & Hashing" sheet (sheet1.xml)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Arrays & Hashing")

# Column A holds dates formatted/stored as plain text (e.g. "08/14/2025"),
# not real date serials - force text formatting first so Excel doesn't
# auto-convert the string into a date value.
$ws.Range("A13:A14").NumberFormat = "@"

# Row 13: Triple ZigZag
$ws.Cells.Item(13, 1).Value = "08/14/2025"
$ws.Cells.Item(13, 2).Value = "Triple ZigZag"
$ws.Cells.Item(13, 3).Value = "Arrays"
$ws.Cells.Item(13, 4).Value = "Easy"
$ws.Cells.Item(13, 5).Value = "Yes"
$ws.Cells.Item(13, 6).Value = "Maybe"
$ws.Cells.Item(13, 7).Value = "Yes`r"

# Row 14: Simple Bank System
$ws.Cells.Item(14, 1).Value = "08/14/2025"
$ws.Cells.Item(14, 2).Value = "Simple Bank System"
$ws.Cells.Item(14, 3).Value = "Arrays"
$ws.Cells.Item(14, 4).Value = "Medium"
$ws.Cells.Item(14, 5).Value = "Yes"
$ws.Cells.Item(14, 6).Value = "Maybe"
$ws.Cells.Item(14, 7).Value = "Yes"

# ---------------------------------------------------------------------------
# 2) Add a new "Others" sheet at the end with the standard header row
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Others"

$headers = @("Date Solved", "Name", "Algorithm", "Difficulty", "Solved First Time", "Revisit?", "Understand?", "Revisit Date #1", "Revisit Date #2", "Revisit Date #3", "Confidence Now")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
